$d = $word.ActiveDocument

# The cover/title-page footer carries the copyright year and the
# "Course Rev" revision string.  Walk every section's primary footer
# (the only one populated in this template) and apply the two textual
# updates described in the commit:
#   1. Copyright year 2016 -> 2017
#   2. "Course Rev 1.0" + ".1"  ->  "Course Rev 1.1.0"
foreach ($sec in $d.Sections) {
    $footer = $sec.Footers(1)
    if ($footer.Exists) {

        $rng = $footer.Range.Duplicate
        $rng.Find.Execute("2016", $true, $false, $false, $false, $false, `
                           $true, 1, $false, "2017", 2)

        $rng = $footer.Range.Duplicate
        $rng.Find.Execute("Course Rev 1.0.1", $true, $false, $false, $false, $false, `
                           $true, 1, $false, "Course Rev 1.1.0", 2)
    }
}
